$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.954.03'
$ws.Range('D3').Value = '1.640.36'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.53'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.16%  '
$ws.Range('E9').Value = '  -1.54%  '
$ws.Range('E10').Value = '  +0.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0882'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.33%  '
$ws.Range('D12').Value = '1.872.88'
$ws.Range('E12').Value = '  +0.35%  '
$ws.Range('D13').Value = '1.636.99'
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('E14').Value = '  +1.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.572'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.89'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.12%  '
$ws.Range('D17').Value = '27.956.74'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '233.33'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('E19').Value = '  +0.52%  '
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('E22').Value = '  +0.58%  '
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('E24').Value = '  -2.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.17'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E26').Value = '  +1.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.49%  '
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('E33').Value = '  +0.93%  '
$ws.Range('D34').Value = '1.416.86'
$ws.Range('E34').Value = '  -3.81%  '
$ws.Range('E35').Value = '  +1.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.35'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.32%  '
$ws.Range('E37').Value = '  +1.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.883'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('E39').Value = '  -0.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.907'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.30%  '
$ws.Range('E41').Value = '  +0.85%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.87'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '66.43'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.17%  '
$ws.Range('E45').Value = '  +3.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.21'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').Value = '1.781.89'
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.86'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.17%  '
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('E50').Value = '  +0.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.61'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.33%  '
